$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.524.43"
$ws.Range("E2").Value = "  +1.27%  "

$ws.Range("D3").Value = "3.154.92"
$ws.Range("E3").Value = "  +3.83%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").Value = "564.56"
$ws.Range("E5").Value = "  +3.04%  "

$ws.Range("D6").Value = "140.98"
$ws.Range("E6").Value = "  +3.17%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").Value = "3.144.56"
$ws.Range("E8").Value = "  +4.07%  "

$ws.Range("E9").Value = "  +2.52%  "

$ws.Range("D10").Value = "6.82"
$ws.Range("E10").Value = "  +6.49%  "

$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +2.98%  "

$ws.Range("E12").Value = "  +3.30%  "

$ws.Range("D13").Value = "36.32"
$ws.Range("E13").Value = "  +3.43%  "

$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("D15").Value = "3.664.29"
$ws.Range("E15").Value = "  +3.77%  "

$ws.Range("D16").Value = "64.599.22"
$ws.Range("E16").Value = "  +0.89%  "

$ws.Range("D17").Value = "3.153.71"
$ws.Range("E17").Value = "  +2.48%  "

$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("D19").Value = "513.73"
$ws.Range("E19").Value = "  +6.70%  "

$ws.Range("E20").Value = "  +4.23%  "

$ws.Range("D21").Value = "13.99"
$ws.Range("E21").Value = "  +4.04%  "

$ws.Range("D22").Value = "0.717"
$ws.Range("E22").Value = "  +5.94%  "

$ws.Range("E23").Value = "  +4.02%  "

$ws.Range("D24").Value = "12.74"
$ws.Range("E24").Value = "  +4.19%  "

$ws.Range("D25").Value = "79.01"
$ws.Range("E25").Value = "  +2.24%  "

$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "8.64"
$ws.Range("E27").Value = "  +10.40%  "

$ws.Range("E28").Value = "  +4.98%  "

$ws.Range("E29").Value = "  +2.92%  "

$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").Value = "26.67"
$ws.Range("E31").Value = "  +4.54%  "

$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("D34").Value = "556.91"
$ws.Range("E34").Value = "  -3.59%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "54.84"
$ws.Range("E35").Value = "  +5.61%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  +3.02%  "

$ws.Range("D37").Value = "5.33"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "0.0427"
$ws.Range("E38").Value = "  +7.12%  "

$ws.Range("E39").Value = "  +4.67%  "

$ws.Range("D40").Value = "3.123.69"
$ws.Range("E40").Value = "  +6.37%  "

$ws.Range("E41").Value = "  +3.07%  "

$ws.Range("E42").Value = "  +1.65%  "

$ws.Range("D43").Value = "2.71"
$ws.Range("E43").Value = "  -6.17%  "

$ws.Range("E44").Value = "  +10.02%  "

$ws.Range("E45").Value = "  +5.23%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "122.61"
$ws.Range("E47").Value = "  +2.48%  "

$ws.Range("D48").Value = "24.96"
$ws.Range("E48").Value = "  +1.98%  "

$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("E51").Value = "  +2.49%  "
